# Update "想去人数" (F) / "最低票价" (G) figures on the 展览 and 全部类型
# sheets to match the freshly re-scraped counts (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 1052
    $ws.Range("F3").Value = 742
    $ws.Range("F6").Value = 1096
    $ws.Range("F8").Value = 1676
    $ws.Range("G8").Value = 55
    $ws.Range("F9").Value = 6150
    $ws.Range("F11").Value = 363
    $ws.Range("F12").Value = 290
    $ws.Range("F13").Value = 91
    $ws.Range("F14").Value = 366
    $ws.Range("F17").Value = 268
    $ws.Range("F18").Value = 1273
    $ws.Range("F20").Value = 114
    $ws.Range("F23").Value = 265
    $ws.Range("F27").Value = 94
    $ws.Range("F28").Value = 3
    $ws.Range("F29").Value = 388
    $ws.Range("F30").Value = 80
    $ws.Range("F32").Value = 77

    if ($sheetName -eq "展览") {
        $ws.Range("F16").Value = 5746
    } else {
        $ws.Range("F16").Value = 5748
    }
}
